$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Set the transaction cost values for rows 5-7 (previously blank) so each
# event's transaction amount is populated; all downstream formulas
# (Numeric/Balance/Check/View sections) recalc automatically.
$ws.Range("B5").Value = 97
$ws.Range("B6").Value = 25
$ws.Range("B7").Value = 25

# Update the active selection to B6
$ws.Range("B6").Select()
